$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values are free-form text that look numeric (prices like "69.662.35",
# "1.00", "0.0000286"). Force the cell to Text format BEFORE writing the value so
# Excel does not silently coerce it into a Double (which would mangle formatting,
# drop trailing zeros, or render tiny values in scientific notation).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.662.35"
$ws.Range("E2").Value = "  -1.87%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.482.57"
$ws.Range("E3").Value = "  -4.31%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.70"
$ws.Range("E5").Value = "  -4.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "192.83"
$ws.Range("E6").Value = "  -3.54%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.611"
$ws.Range("E7").Value = "  -2.86%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.470.83"
$ws.Range("E8").Value = "  -4.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.206"
$ws.Range("E10").Value = "  -7.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.619"
$ws.Range("E11").Value = "  -4.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "51.41"
$ws.Range("E12").Value = "  -5.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000286"
$ws.Range("E13").Value = "  -7.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.15"
$ws.Range("E14").Value = "  -4.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.040.05"
$ws.Range("E15").Value = "  -4.06%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "646.31"
$ws.Range("E16").Value = "  -2.43%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.498.82"
$ws.Range("E17").Value = "  -2.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.477.21"
$ws.Range("E18").Value = "  -4.68%  "
$ws.Range("E19").Value = "  -5.33%  "
$ws.Range("E20").Value = "  -1.83%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.20"
$ws.Range("E21").Value = "  -4.79%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.947"
$ws.Range("E22").Value = "  -5.42%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.83"
$ws.Range("E23").Value = "  -3.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.28"
$ws.Range("E24").Value = "  -0.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "99.16"
$ws.Range("E25").Value = "  -5.81%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.29"
$ws.Range("E26").Value = "  -7.61%  "
$ws.Range("E27").Value = "  -4.45%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.97"
$ws.Range("E28").Value = "  -5.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.34"
$ws.Range("E29").Value = "  -4.76%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.53"
$ws.Range("E30").Value = "  -4.94%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.32"
$ws.Range("E31").Value = "  -7.64%  "
$ws.Range("E32").Value = "  -6.98%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.64"
$ws.Range("E33").Value = "  -4.98%  "
$ws.Range("E34").Value = "  -5.54%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "60.96"
$ws.Range("E35").Value = "  -4.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.750.43"
$ws.Range("E36").Value = "  -6.17%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "529.41"
$ws.Range("E37").Value = "  +3.73%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0₃0792"
$ws.Range("E39").Value = "  -9.94%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.94"
$ws.Range("E40").Value = "  -3.35%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.51"
$ws.Range("E41").Value = "  -1.22%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.373"
$ws.Range("E42").Value = "  -4.39%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.134"
$ws.Range("E43").Value = "  -1.50%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "34.35"
$ws.Range("E44").Value = "  -6.91%  "
$ws.Range("B45").Value = "CoreDAO"
$ws.Range("C45").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.49"
$ws.Range("E45").Value = "  +69.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0442"
$ws.Range("E46").Value = "  -4.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.36"
$ws.Range("E47").Value = "  -4.54%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.83"
$ws.Range("E48").Value = "  -8.42%  "
$ws.Range("E49").Value = "  -4.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  -0.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.18"
$ws.Range("E51").Value = "  -6.16%  "
